# "Fixed product name typo & adding tiny images"
#
# 1. Corrects the "Butterfly Earings" -> "Butterfly Earrings" typo for the
#    SAMPLE006 product (row 5 on the "Main" sheet, and the mirrored SLUG
#    entry on the "Categories" sheet).
# 2. Swaps the placeholder "indiana-jones-hat.jpg" image filename used by
#    the Brown Fedora product (row 4) for its own "brown-fedora-01.jpg".

$wb = $excel.ActiveWorkbook

$main = $wb.Worksheets.Item("Main")
$categories = $wb.Worksheets.Item("Categories")

# --- Butterfly Earrings (row 5): fix the spelling everywhere --------------
$main.Range("A5").Value = "butterfly-earrings"
$main.Range("E5").Value = "Butterfly Earrings"
$main.Range("M5").Value = "Sample Butterfly Earrings Lorem ipsum dolor sit amet, consectetur adipisicing elit, sed do eiusmod tempor incididunt ut labore et dolore magna aliqua. Ut enim ad minim veniam, quis nostrud exercitation ullamco laboris nisi ut aliquip ex ea commodo consequat. Duis aute irure dolor in reprehenderit in voluptate velit esse cillum dolore eu fugiat nulla pariatur. Excepteur sint occaecat cupidatat non proident, sunt in culpa qui officia deserunt mollit anim id est laborum"

# The "Categories" sheet keeps its own copy of each product's SLUG.
$categories.Range("A4").Value = "butterfly-earrings"

# --- Brown Fedora (row 4): give it its own product image -----------------
$main.Range("L4").Value = "brown-fedora-01.jpg"

# --- Restore the on-screen selection/scroll position seen in the edit -----
$categories.Range("A9").Select()

$main.Activate()
$main.Range("L5").Select()
